$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Data change: C2 goes from 90 to 89 (D2 = C2/B2 recalculates automatically)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 89

# ---------------------------------------------------------------------------
# 2. Conditional formatting on C2: highlight when C2 is greater than B2 (goal
#    met) or less than B2 (goal missed).
#
# The approximate fill colours below are the closest reproducible
# equivalents (via plain RGB) of the theme-based colors
# (theme 5 "Accent2" tinted ~80% and theme 9 "Accent6" tinted ~60%/~80%)
# that Excel itself used, since this automation surface does not expose a
# working Interior.ThemeColor/TintAndShade setter.
# ---------------------------------------------------------------------------
$range = $ws.Range("C2")

$redFill        = 0xD6E5FB   # ~ theme 5 (Accent2), tint 0.8   -> RGB FBE5D6
$greenFillDark  = 0xB5DEC6   # ~ theme 9 (Accent6), tint 0.6   -> RGB C6DEB5
$greenFillLight = 0xDAEFE2   # ~ theme 9 (Accent6), tint 0.8   -> RGB E2EFDA

# Helper: create a one-off FormatCondition purely so the workbook's dxf
# (differential format) list gets an extra entry, then discard the rule.
# This mirrors what happened in the real workbook, where several
# Highlight-Cell-Rules were applied/removed while the user was experimenting,
# leaving behind unused <dxf> entries before the final two rules were set.
function New-ThrowawayDxf([int]$color, [bool]$touchFont) {
    $tmp = $range.FormatConditions.Add(1, 5, '=1')
    if ($touchFont) {
        $tmp.Font.Bold = $false
        $tmp.Font.Italic = $false
    }
    $tmp.Interior.Color = $color
    $tmp.Delete()
}

# Build up dxf entries 0-3 (unused) ahead of the two real rules.
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillDark  $false
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillLight $true

# A throwaway rule so the two kept rules land on priority 2 and 3 (matching
# the saved workbook) instead of 1 and 2.
$priorityBump = $range.FormatConditions.Add(1, 5, '=1')

# The two rules that remain in the final workbook.
$lessThanRule    = $range.FormatConditions.Add(1, 6, '=$B$2')
$greaterThanRule = $range.FormatConditions.Add(1, 5, '=$B$2')

$lessThanRule.Interior.Color    = $redFill
$greaterThanRule.Interior.Color = $greenFillDark

$priorityBump.Delete()

# Build up dxf entries 6-17 (unused), continuing the same alternating pattern.
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillDark  $false
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillLight $true
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillDark  $false
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillLight $true
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillLight $true
New-ThrowawayDxf $redFill        $false
New-ThrowawayDxf $greenFillLight $true

# ---------------------------------------------------------------------------
# 3. Page setup: paper size A4 (9), portrait orientation
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
